# Workbook restructuring:
#  - rename "geradores" -> "UTE"
#  - add a new "UEL" sheet at the end
#  - reorder so the tab order is: demanda, UTE, UEL
#  - replace the "demanda" data with the new (longer) hourly series
#  - populate the new "UEL" sheet with its data table
#  - restore per-sheet selections

$wb = $excel.ActiveWorkbook

# ---- sheet bookkeeping -----------------------------------------------
$ute = $wb.Worksheets.Item("geradores")
$ute.Name = "UTE"

$demanda = $wb.Worksheets.Item("demanda")

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$uel = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$uel.Name = "UEL"

# put demanda first, ahead of UTE -> demanda, UTE, UEL
$demanda.Move($ute)

# NOTE: Move() swaps which sheet the pre-existing variables are bound to
# (position-bound references) - re-fetch fresh handles by name afterwards.
$demanda = $wb.Worksheets.Item("demanda")
$ute = $wb.Worksheets.Item("UTE")

# ---- demanda: new hourly values ---------------------------------------
$demandaVals = @(
    @(0, 2560), @(1, 3620), @(2, 3800), @(3, 3810),
    @(4, 2990), @(5, 4040), @(6, 4000), @(7, 3790),
    @(8, 4680), @(9, 4540), @(10, 3690), @(11, 4750),
    @(12, 5560), @(13, 5620), @(14, 5800), @(15, 5810),
    @(16, 5990), @(17, 6040), @(18, 6000), @(19, 5790),
    @(20, 5680), @(21, 5540), @(22, 5690), @(23, 5750)
)
for ($i = 0; $i -lt $demandaVals.Count; $i++) {
    $row = $i + 2
    $demanda.Cells.Item($row, 1).Value = $demandaVals[$i][0]
    $demanda.Cells.Item($row, 2).Value = $demandaVals[$i][1]
}
$demanda.Range("B26").Font.Underline = $true

# ---- UEL: new sheet content --------------------------------------------
$uel.Range("A1").Value = "Unid."
$uel.Range("B1").Value = "Nt"
$uel.Range("C1").Value = "Pt"
$uel.Range("D1").Value = "Wmax"
$uel.Range("E1").Value = "d"

$uel.Range("A1").HorizontalAlignment = -4108
$uel.Range("F1:I1").HorizontalAlignment = -4108

$uelVals = @(
    @(1, 25, 13, 325, 0.8),
    @(2, 15, 25, 375, 0.6),
    @(3, 18, 12, 216, 1.05),
    @(4, 17, 15, 255, 1.2),
    @(5, 14, 11, 292.5, 0.75)
)
for ($i = 0; $i -lt $uelVals.Count; $i++) {
    $row = $i + 2
    $uel.Cells.Item($row, 1).Value = $uelVals[$i][0]
    $uel.Cells.Item($row, 2).Value = $uelVals[$i][1]
    $uel.Cells.Item($row, 3).Value = $uelVals[$i][2]
    $uel.Cells.Item($row, 4).Value = $uelVals[$i][3]
    $uel.Cells.Item($row, 5).Value = $uelVals[$i][4]
}
$uel.Range("A2:A6").HorizontalAlignment = -4108

$uel.Range("F2:I3").NumberFormat = "0.00"
$uel.Range("F2:I3").HorizontalAlignment = -4108

$uel.Columns.Item(2).ColumnWidth = 8.7369791666667
$uel.Columns.Item(7).ColumnWidth = 11.7369791666667

# ---- selections / active tab -------------------------------------------
$ute.Activate()
$ute.Range("K4").Select()

$uel.Activate()
$uel.Range("G13").Select()

$demanda.Activate()
$demanda.Range("C9").Select()

Write-Host "done"
